# Slide 11 ("Sources Used:") / Shape 2 ("Content Placeholder 2") lists a set
# of hyperlinked source URLs, one per paragraph.
#
# This edit:
#  1. Splits the NexusMods "https://www.nexusmods.com/newvegas/mods/68714"
#     paragraph's single run into two runs (same rId7 hyperlink):
#       "https://www.nexusmods.com/newvegas/mods" + "/68714"
#  2. Removes the following "knowyourmeme" paragraph entirely, together with
#     the blank paragraph right after it, so the NexusMods paragraph merges
#     directly with what used to be the second blank paragraph (inheriting
#     its plain end-of-paragraph mark).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# --- Step 1: split the "mods/68714" run into two runs ----------------------
$para6 = $tr.Paragraphs(6, 1)
$linkStart = $para6.Start
$prefix = "https://www.nexusmods.com/newvegas/mods"
$suffix = "/68714"

$tailPart = $tr.Characters($linkStart + $prefix.Length, $suffix.Length)
$tailPart.Text = $suffix

# --- Step 2: delete the "knowyourmeme" paragraph + the blank one after it --
$para7 = $tr.Paragraphs(7, 1)
$knowYourMeme = "https://knowyourmeme.com/photos/1809295-patrick-star"
$charsToDelete = $knowYourMeme.Length + 1 + 1   # its text + its own para mark + the next (blank) para mark
$toRemove = $tr.Characters($para7.Start, $charsToDelete)
$toRemove.Delete()
